$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(16, "1047463924", "BRANDO HERRERA VERGARA", "1710", 29509),
    @(17, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2010", 26041),
    @(18, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2009", 31249),
    @(19, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2008", 31249),
    @(20, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2007", 31249),
    @(21, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2006", 31249),
    @(22, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2005", 31249),
    @(23, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2004", 31249),
    @(24, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2003", 31249),
    @(25, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2002", 31249),
    @(26, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "2001", 31249),
    @(27, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1912", 31249),
    @(28, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1911", 31249),
    @(29, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1910", 31249),
    @(30, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1909", 31249),
    @(31, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1908", 31249),
    @(32, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1907", 31249),
    @(33, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1906", 31249),
    @(34, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1905", 31249),
    @(35, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1904", 31249),
    @(36, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1903", 31249),
    @(37, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1902", 31249),
    @(38, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1901", 31249),
    @(39, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1812", 31249),
    @(40, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1811", 31249),
    @(41, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1810", 31249),
    @(42, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1809", 31249),
    @(43, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1808", 29509),
    @(44, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1807", 29509),
    @(45, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1806", 29509),
    @(46, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1805", 29509),
    @(47, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1804", 29509),
    @(48, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1803", 29509),
    @(49, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1802", 29509),
    @(50, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1801", 29509),
    @(51, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1712", 29509),
    @(52, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1711", 29509),
    @(53, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1710", 29509),
    @(54, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1709", 29509),
    @(55, "1002190015", "DEYBIS ALEXANDER YOUNG ACEVEDO", "1708", 29509),
    @(56, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2010", 26041),
    @(57, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2009", 31249),
    @(58, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2008", 31249),
    @(59, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2007", 31249),
    @(60, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2006", 31249),
    @(61, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2005", 31249),
    @(62, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2004", 31249),
    @(63, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2003", 31249),
    @(64, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2002", 31249),
    @(65, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "2001", 31249),
    @(66, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1912", 31249),
    @(67, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1911", 31249),
    @(68, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1910", 31249),
    @(69, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1909", 31249),
    @(70, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1908", 31249),
    @(71, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1907", 31249),
    @(72, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1906", 31249),
    @(73, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1905", 31249),
    @(74, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1904", 31249),
    @(75, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1903", 31249),
    @(76, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1902", 31249),
    @(77, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1901", 31249),
    @(78, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1812", 31249),
    @(79, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1811", 31249),
    @(80, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1810", 31249),
    @(81, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1809", 31249),
    @(82, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1808", 29509),
    @(83, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1807", 29509),
    @(84, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1806", 29509),
    @(85, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1805", 29509),
    @(86, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1804", 29509),
    @(87, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1803", 29509),
    @(88, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1802", 29509),
    @(89, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1801", 29509),
    @(90, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1712", 29509),
    @(91, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1711", 29509),
    @(92, "1047484535", "MANUEL ENRIQUE DORIA GOMEZ", "1710", 29509)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}
